# Appends 25 more blood-glucose log entries (rows 516-540) to Sheet1,
# continuing the existing 5-minute-interval time series in column A
# (last existing reading: 2026/02/13 08:58) with matching values in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 516

$times = @(
    "2026/02/13 09:03",
    "2026/02/13 09:08",
    "2026/02/13 09:13",
    "2026/02/13 09:18",
    "2026/02/13 09:23",
    "2026/02/13 09:28",
    "2026/02/13 09:33",
    "2026/02/13 09:38",
    "2026/02/13 09:43",
    "2026/02/13 09:48",
    "2026/02/13 09:53",
    "2026/02/13 09:58",
    "2026/02/13 10:03",
    "2026/02/13 10:08",
    "2026/02/13 10:13",
    "2026/02/13 10:18",
    "2026/02/13 10:23",
    "2026/02/13 10:28",
    "2026/02/13 10:33",
    "2026/02/13 10:38",
    "2026/02/13 10:43",
    "2026/02/13 10:48",
    "2026/02/13 10:53",
    "2026/02/13 10:58",
    "2026/02/13 11:03"
)

$readings = @(
    "13.5",
    "13.6",
    "14.0",
    "14.1",
    "13.4",
    "14.0",
    "14.1",
    "14.3",
    "14.5",
    "13.9",
    "13.4",
    "13.5",
    "13.8",
    "13.8",
    "13.7",
    "13.7",
    "13.4",
    "13.3",
    "13.5",
    "13.7",
    "13.8",
    "13.4",
    "13.7",
    "13.5",
    "12.5"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $r = $startRow + $i
    $ws.Range("A$r").Value = $times[$i]

    # The reading looks numeric ("13.5"), but the workbook stores column B
    # as plain text (shared strings) with no special number format on the
    # cells. Assigning it straight to .Value would make Excel infer a
    # Number type, and forcing text via NumberFormat="@" (or a leading
    # apostrophe) would permanently stamp a new cell style ("@" / quote
    # prefix) onto the cell - something the target workbook does not have.
    # Instead, write it as a formula that evaluates to a text string, then
    # immediately paste-special (values only) over itself: this commits a
    # literal text cell identical in shape to the existing data, with zero
    # style-table side effects.
    $ws.Range("B$r").Formula = '="' + $readings[$i] + '"'
}

$dataRange = $ws.Range("B" + $startRow + ":B" + ($startRow + $times.Length - 1))
$dataRange.Copy()
$dataRange.PasteSpecial(-4163)
